$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.475.36'
$ws.Range('E2').Value = '  +5.76%  '
$ws.Range('D3').Value = '2.053.16'
$ws.Range('E3').Value = '  +4.18%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  +3.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.651'
$ws.Range('E6').Value = '  +2.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.20'
$ws.Range('E7').Value = '  +17.05%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +7.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.70'
$ws.Range('E10').Value = '  +2.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0771'
$ws.Range('E11').Value = '  +5.64%  '
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.912'
$ws.Range('E13').Value = '  -2.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.91'
$ws.Range('E14').Value = '  +4.85%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.54'
$ws.Range('E15').Value = '  +29.25%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.352.80'
$ws.Range('E16').Value = '  +4.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.59'
$ws.Range('E17').Value = '  +6.99%  '
$ws.Range('D18').Value = '2.060.89'
$ws.Range('E18').Value = '  +4.19%  '
$ws.Range('D19').Value = '37.268.34'
$ws.Range('E19').Value = '  +5.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.82'
$ws.Range('E20').Value = '  +3.56%  '
$ws.Range('D21').Value = '0.0₃0881'
$ws.Range('E21').Value = '  +5.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.45'
$ws.Range('E22').Value = '  +6.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.60'
$ws.Range('E23').Value = '  +3.87%  '
$ws.Range('E24').Value = '  +6.26%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  +5.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.95'
$ws.Range('E27').Value = '  +10.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.19'
$ws.Range('E28').Value = '  -0.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.03'
$ws.Range('E29').Value = '  +5.39%  '
$ws.Range('E30').Value = '  +28.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.33'
$ws.Range('E31').Value = '  +10.61%  '
$ws.Range('E32').Value = '  +3.94%  '
$ws.Range('E33').Value = '  +10.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.73'
$ws.Range('E34').Value = '  +10.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0626'
$ws.Range('E35').Value = '  +6.95%  '
$ws.Range('E36').Value = '  +4.80%  '
$ws.Range('E37').Value = '  +5.38%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('E39').Value = '  +17.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.00'
$ws.Range('E40').Value = '  +35.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.104'
$ws.Range('E41').Value = '  +17.70%  '
$ws.Range('E42').Value = '  +3.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.02'
$ws.Range('E43').Value = '  +5.20%  '
$ws.Range('E44').Value = '  +6.87%  '
$ws.Range('E45').Value = '  +4.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.17'
$ws.Range('E46').Value = '  +9.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '96.92'
$ws.Range('E47').Value = '  +7.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.01'
$ws.Range('E48').Value = '  +7.60%  '
$ws.Range('D49').Value = '1.427.47'
$ws.Range('E49').Value = '  +4.61%  '
$ws.Range('E50').Value = '  +2.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.96'
$ws.Range('E51').Value = '  +2.13%  '
